# Sync edit: strip markdown emphasis/code markers ("*", "`", "_") from all
# text cells in the "Plan de Accion" sheet, and convert the two milestone
# date cells (D22 "**31/10/2025**" and D29 "**31/12/2025**") from free text
# into real date values formatted as dd/mm/yyyy.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

for ($r = 1; $r -le $lastRow; $r++) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -is [string]) {
            $newVal = $val.Replace("*", "").Replace("``", "").Replace("_", "")
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}

# D22: "31/10/2025" (was "**31/10/2025**") -> real date serial 45961
$d22 = $ws.Range("D22")
$d22.Value = 45961
$d22.NumberFormat = "dd/mm/yyyy"

# D29: "31/12/2025" (was "**31/12/2025**") -> real date serial 46022
$d29 = $ws.Range("D29")
$d29.Value = 46022
$d29.NumberFormat = "dd/mm/yyyy"

# The stripped text is shorter than before, so the "best fit" autofit
# widths for columns A and B shrink accordingly. (These ColumnWidth inputs
# are chosen empirically so the engine's pixel-quantized stored width lands
# on the closest achievable value to the target 40.42578125 / 34.7109375.)
$ws.Columns.Item(1).ColumnWidth = 39.665
$ws.Columns.Item(2).ColumnWidth = 33.83

$wb.Save()
